$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two Salesforce record-id placeholders in column D (rows 2 and 3) were
# refreshed to new ids. Write D3 first so the shared-strings table gains the
# new unique strings in the same order Excel would (row 3's new id before
# row 2's), matching the row-major save order.
$ws.Range("D3").Value = "a0Nq0000003PKUh"
$ws.Range("D2").Value = "a0Nq0000003PKUc"

# Update the sheet's active selection to D2 (previously S1).
$ws.Range("D2").Select() | Out-Null
